# Regenerate save_data to use K (strikeouts) instead of Strike# (total strikes)
# in column G, for rows 2-36 (the per-game pitching log rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts for each game row, replacing the old Strike# totals.
$kValues = @{
    2  = 2
    3  = 3
    4  = 2
    5  = 5
    6  = 8
    7  = 7
    8  = 3
    9  = 7
    10 = 5
    11 = 9
    12 = 6
    13 = 2
    14 = 5
    15 = 1
    16 = 7
    17 = 3
    18 = 8
    19 = 6
    20 = 5
    21 = 1
    22 = 7
    23 = 3
    24 = 7
    25 = 6
    26 = 3
    27 = 3
    28 = 4
    29 = 4
    30 = 4
    31 = 8
    32 = 5
    33 = 4
    34 = 3
    35 = 1
    36 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
